$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update barcode_type values from "ean" to "ean13"
$ws.Range("I2").Value = "ean13"
$ws.Range("I3").Value = "ean13"

# Update the active selection to I3
$ws.Range("I3").Select()
